$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add the new "url" column header.
$ws.Range("D1").Value = "url"

# 2. Insert a new row for "Navy Moves" in its alphabetically-sorted position
#    (between "Mad Mix Game" in row 41 and "Oh Mummy", which was row 42).
$ws.Rows.Item(42).Insert()
$ws.Range("A42").Value = "Navy Moves"
$ws.Range("B42").Value = 1988
$ws.Range("C42").Value = "Dinamic Software"

# 3. Add a hyperlink on the "La Abadía del Crimen" row (row 39) pointing to
#    its source page.
$target = $ws.Range("D39")
$ws.Hyperlinks.Add($target, "https://amstrad.es/doku.php?id=juegos:la_abadia_del_crimen") | Out-Null
